# TASK_35 — robot game: merge 5.4 "check winner" and 5.5 "close game" into a
# single backlog row, mark 5.3 ("robot brain") as finished, and move the
# "active tab" from the backlog sheet to the tasks sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Задачи"           (tasks)
$ws2 = $wb.Worksheets.Item(2)   # "Бэклог задач"      (backlog)

# ---------------------------------------------------------------------
# 1) Sheet "Задачи": row 37 (task #34, "5.3 – Мозг робота") now has a
#    finish date in column E (it used to be empty / still open).
# ---------------------------------------------------------------------
$ws1.Range("E37").Value = 42014.952777777777

# ---------------------------------------------------------------------
# 2) Row 38 (task #35) becomes the merged "5.4 / 5.5" task. Pick up the
#    border formatting already used elsewhere on the sheet (A6/C6 carry
#    the exact bordered xf we need) instead of re-deriving borders by
#    hand, so we land on the same shared cellXf the workbook already
#    uses elsewhere.
# ---------------------------------------------------------------------
$ws1.Range("A6").Copy() | Out-Null
$ws1.Range("A38").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> bordered style (no numfmt)

$ws1.Range("A6").Copy() | Out-Null
$ws1.Range("B38").PasteSpecial(-4122) | Out-Null   # bordered style, then add wrap text
$ws1.Range("B38").WrapText = $true

$ws1.Range("C6").Copy() | Out-Null
$ws1.Range("C38").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> bordered + date numfmt
$ws1.Range("D38").PasteSpecial(-4122) | Out-Null
$ws1.Range("E38").PasteSpecial(-4122) | Out-Null

$ws1.Range("B38").Value = "5.4 – Проверка победителя. (LogicRobot.checkWinner)" + [char]10 + "5.5 – Покинуть и игру с роботом. (LogicRobot.closeGame)"
$ws1.Range("C38").Value = 42013.770138888889
$ws1.Range("D38").Value = 42014.952777777777
$ws1.Range("E38").Value = 42015.051388888889

$ws1.Rows.Item(38).RowHeight = 30

# ---------------------------------------------------------------------
# 3) Row 39 (former task #36, "5.5 – Покинуть игру") is now absorbed
#    into row 38, so its number/name/created-date go away; D/E stay the
#    same untouched blanks they always were.
# ---------------------------------------------------------------------
$ws1.Range("A39").ClearContents()
$ws1.Range("B39").Clear()
$ws1.Range("C39").Clear()

# ---------------------------------------------------------------------
# 4) View state: the tasks sheet becomes the active tab/selection
#    (previously it was the backlog sheet). Touch the backlog sheet's
#    selection first so the very last Select() call is what leaves the
#    tasks sheet as the active one.
# ---------------------------------------------------------------------
$ws2.Range("C25").Select() | Out-Null
$ws1.Range("C49").Select() | Out-Null
